$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '29.978.95'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.03%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.116.99'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +0.81%  '

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.011'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.69%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '347.10'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.24%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.009'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.62%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5199'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.94%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.4451'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.68%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '53.88'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +2.88%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.09334'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.61%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.181'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.95%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '25.16'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -0.27%  '

# Row 13
$ws.Range("B13").Value = 'Chainlink'
$ws.Range("C13").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '8.556'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +4.96%  '

# Row 14
$ws.Range("B14").Value = 'WrappedEther'
$ws.Range("C14").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '2.133.65'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +1.64%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.904'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +2.41%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '102.81'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +3.36%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.00001161'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.02%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '1.010'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.67%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '21.57'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +4.91%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.06701'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.27%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.311'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +1.47%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '1.009'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.65%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '30.039.14'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.02%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '12.71'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.72%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.325'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.52%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.392.72'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +2.01%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '22.12'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.50%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.539'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.54%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '162.80'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.10%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '134.13'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.64%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.151'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -1.28%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.773'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +8.23%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.1056'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.25%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '6.247'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.26%  '

# Row 35
$ws.Range("B35").Value = 'HuobiToken'
$ws.Range("C35").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.979'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.72%  '

# Row 36
$ws.Range("B36").Value = 'InternetComputer(DFINITY)'
$ws.Range("C36").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '6.545'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +5.58%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '10.73'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +5.65%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.02619'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +2.20%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.06866'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +1.32%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.7073'
$ws.Range("D40").Style = "Normal"

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '12.69'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +1.22%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.2241'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -1.59%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.327'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +1.15%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.6844'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +2.73%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '14.56'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +2.83%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.369'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +4.10%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.009'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.69%  '

# Row 48
$ws.Range("B48").Value = 'BabyDogeCoin'
$ws.Range("C48").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.00000000364'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +2.15%  '

# Row 49
$ws.Range("B49").Value = 'WEMIXTOKEN'
$ws.Range("C49").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.278'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +10.04%  '

# Row 50
$ws.Range("B50").Value = 'PancakeSwap'
$ws.Range("C50").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '3.633'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.09%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.224'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.30%  '
